# "Test 20 Works, Final Touches Made"
# Fill in the final couple of missing timesheet entries for week 15 and
# start week 16, then leave the selection parked on the next cell to fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timelog")

# Week 15 (row 32 = Saturday/Sunday hours, row 33 = running totals):
#   Saturday (H32) bumped from 1 to 2.5 hours.
#   Sunday (I32) filled in at 5.5 hours (was blank).
$ws.Range("H32").Value = 2.5
$ws.Range("I32").Value = 5.5

# Week 16 (row 34): Monday (C34) logged at 2 hours (was blank).
$ws.Range("C34").Value = 2

# Move the active selection to I34, the next day to be filled in.
$ws.Range("I34").Select()
